$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.143.33'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '1.678.25'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  -0.52%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.66'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5298'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.19%  '
$ws.Range('E7').Value = '  -0.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2686'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06329'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.34'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07543'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').Value = '1.689.15'
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.518'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5693'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008177'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.85'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').Value = '26.177.88'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.872'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.60'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.07'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.226'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.004'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '148.66'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1263'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.667'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.09'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06333'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.344'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.289'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.553'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.549'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.680'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.013'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6083'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.414'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.722'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.168'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01619'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.27%  '
$ws.Range('D40').Value = '1.096.52'
$ws.Range('E40').Value = '  -2.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8725'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.12'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('D44').Value = '1.829.20'
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('E45').Value = '  -1.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.19'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.009'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.012'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.57%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05260'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4266'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.978'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.00%  '
